# Generate Report for Handoff
# The "b.md" row (row 3) on each sheet moves from "handed back" state to
# "ready for handoff" state: status text, handoff file name and handoff
# timestamp are refreshed to reflect the new handoff package
# (b.63290e5768f688058c7b37413b0a5c26c308f864.*).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row 3 (b.md) status + latest handoff date
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-22-11 16:22:37"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 (b.md) status, handoff file + handoff datetime
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-11 16:22:34"

foreach ($hl in $zhcn.Hyperlinks) {
    if ($hl.Range.Address() -eq '$D$3') {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------
# de-de sheet: row 3 (b.md) status, handoff file + handoff datetime
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value = "2016-03-11 16:22:37"

foreach ($hl in $dede.Hyperlinks) {
    if ($hl.Range.Address() -eq '$D$3') {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
